# Textbox response formatting fix
# Updates task-order sheet names and refreshed CSV stim-file / condition
# values (timestamps regenerated on a later run of the task).

$wb = $excel.ActiveWorkbook

# --- Rename sheets (new timestamp suffixes) ---------------------------------
$wb.Worksheets.Item(1).Name = "GNG_TO-16511686920700538"
$wb.Worksheets.Item(2).Name = "NB_TO-165116869326342"
$wb.Worksheets.Item(3).Name = "RS_TO-16511686932643878"
$wb.Worksheets.Item(4).Name = "TOL_TO-16511686933123834"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16511686933743837"

# --- Sheet 1 : GNG_TO --------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-1651168692036018.csv"
$ws1.Range("B3").Value = "GNG_stims-1651168692053049.csv"
$ws1.Range("B4").Value = "go_stims-16511686920540142.csv"
$ws1.Range("B5").Value = "GNG_stims-16511686920690196.csv"

# --- Sheet 2 : NB_TO ----------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "ZB-match_6-1651168692285567.csv"
$ws2.Range("B3").Value = "TB-16511686932423866.csv"
$ws2.Range("B4").Value = "OB-16511686928093953.csv"
$ws2.Range("B5").Value = "TB-16511686930603838.csv"
$ws2.Range("B6").Value = "OB-1651168692650393.csv"
$ws2.Range("B7").Value = "ZB-match_3-1651168692323383.csv"
$ws2.Range("B8").Value = "TB-16511686928783832.csv"
$ws2.Range("B9").Value = "ZB-match_5-16511686923623843.csv"
$ws2.Range("B10").Value = "OB-1651168692700386.csv"

# --- Sheet 3 : RS_TO (eyes closed / eyes open swap) --------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# --- Sheet 4 : TOL_TO ---------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16511686932793841.csv"
$ws4.Range("B3").Value = "ZM_stims-16511686932673867.csv"
$ws4.Range("B4").Value = "MM_stims-1651168693295393.csv"
$ws4.Range("B5").Value = "ZM_stims-1651168693280385.csv"
$ws4.Range("B6").Value = "MM_stims-16511686933113883.csv"
$ws4.Range("B7").Value = "ZM_stims-16511686932963898.csv"

# --- Sheet 5 : vSAT_TO ---------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16511686933274212.csv"
$ws5.Range("B3").Value = "SAT_stims-16511686933153832.csv"
$ws5.Range("B4").Value = "vSAT_stims-16511686933593886.csv"
$ws5.Range("B5").Value = "vSAT_stims-1651168693343421.csv"
